# Fix duplicate rows where the damage columns (W..AB) had been incorrectly
# dumped as a single serialized array string into column B ("sl added in duplicate rows").
# This script distributes the array values into their proper columns and blanks out B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)   # "Stock Report" sheet

$rows = @(7,11,13,15,17,19,21,24,25,28,29,31,33)

foreach ($r in $rows) {
    # Clear the stray serialized-array text previously stored in column B
    $ws.Range("B$r").Value = $null
    $ws.Range("B$r").Style = "Normal"

    # Recreate blank cells for all the normal (non-date) columns on the row
    $ws.Range("C$r`:L$r").Style = "Normal"
    $ws.Range("N$r`:T$r").Style = "Normal"
    $ws.Range("V$r").Style = "Normal"

    # Gate In Date / Di Date columns keep the date number format, but stay blank
    $ws.Range("M$r").NumberFormat = "YYYY-MM-DD"
    $ws.Range("U$r").NumberFormat = "YYYY-MM-DD"
}

$ws.Range("W7").Value = "FLOORS-(F)"
$ws.Range("X7").Value = "FLOOR BOARD-(FLOOR BOARD)"
$ws.Range("Y7").Value = "F/BOARD DIRTY BY BADLY SCRAP DUST & SCRATCHED."
$ws.Range("Z6:AB6").Copy($ws.Range("Z7:AB7"))

$ws.Range("W11").Value = "FLOORS-(F)"
$ws.Range("X11").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y11").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST & ODOUR "
$ws.Range("Z6:AB6").Copy($ws.Range("Z11:AB11"))

$ws.Range("W13").Value = "FLOORS-(F)"
$ws.Range("X13").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y13").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z13:AB13"))

$ws.Range("W15").Value = "FLOORS-(F)"
$ws.Range("X15").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y15").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z15:AB15"))

$ws.Range("W17").Value = "FLOORS-(F)"
$ws.Range("X17").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y17").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z17:AB17"))

$ws.Range("W19").Value = "FLOORS-(F)"
$ws.Range("X19").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y19").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z19:AB19"))

$ws.Range("W21").Value = "FLOORS-(F)"
$ws.Range("X21").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y21").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z21:AB21"))

$ws.Range("W24").Value = "FLOORS-(F)"
$ws.Range("X24").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y24").Value = "FLOOR BOARD SLIGHTLY UP WARD."
$ws.Range("Z6:AB6").Copy($ws.Range("Z24:AB24"))

$ws.Range("W25").Value = "FLOORS-(F)"
$ws.Range("X25").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y25").Value = "FLOOR BOARD DIRTY BY DUST."
$ws.Range("Z6:AB6").Copy($ws.Range("Z25:AB25"))

$ws.Range("W28").Value = "FLOORS-(F)"
$ws.Range("X28").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y28").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z28:AB28"))

$ws.Range("W29").Value = "FLOORS-(F)"
$ws.Range("X29").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y29").Value = "F/B STEEL PLATE FITTING 12'X08' & LOOSED ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z29:AB29"))

$ws.Range("W31").Value = "FLOORS-(F)"
$ws.Range("X31").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y31").Value = "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST ."
$ws.Range("Z6:AB6").Copy($ws.Range("Z31:AB31"))

$ws.Range("W33").Value = "FLOORS-(F)"
$ws.Range("X33").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y33").Value = "F/B DIRTY BY BADLY MUD+SAND,TYER MARK STAPLE FIBER DUST & ODOUR"
$ws.Range("Z6:AB6").Copy($ws.Range("Z33:AB33"))

# Column widths were re-measured after the data moved out of column B into W:AB
$ws.Cells.Item(1,2).ColumnWidth = 7.85715
$ws.Cells.Item(1,2).EntireColumn.Hidden = $true
$ws.Cells.Item(1,25).ColumnWidth = 78.2857
